$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Clear existing hyperlinks so we can rebuild them cleanly in the correct row positions
$ws.Cells.Hyperlinks.Delete()

# Row 2
$ws.Range("A2").Value = '2025-09-26 12:35:31'
$ws.Range("B2").Value = 'LINExChatGPTx美容室向け予約Bot (仕様書、契約書あり)'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5400801'
$ws.Range("G2").Value = 445
$ws.Range("H2").Value = '🔥GPT,ChatGPT ★bot'
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5400801')
$ws.Range("F2").Style = "Hyperlink"

# Row 3
$ws.Range("A3").Value = '2025-09-26 12:35:31'
$ws.Range("B3").Value = '自社開発のロジシステムをサポート及び開発できる方募集【PHP, Python, VBA etc】'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E3").Value = '期限情報なし'
$ws.Range("F3").Value = 'https://www.lancers.jp/work/detail/5389460'
$ws.Range("G3").Value = 305
$ws.Range("H3").Value = '🔥Python ◆開発 ○PHP'
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5389460')
$ws.Range("F3").Style = "Hyperlink"

# Row 4
$ws.Range("A4").Value = '2025-09-26 12:35:31'
$ws.Range("B4").Value = '【急募】WindwosサーバーでのDjangoアプリ環境構築依頼'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5401279'
$ws.Range("G4").Value = 148
$ws.Range("H4").Value = '🔥Django ◇アプリ'
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5401279')
$ws.Range("F4").Style = "Hyperlink"

# Row 5
$ws.Range("A5").Value = '2025-09-26 12:35:31'
$ws.Range("B5").Value = '【開発依頼】Amazonセラー向け 価格自動調整ツールの開発'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("F5").Value = 'https://www.lancers.jp/work/detail/5401202'
$ws.Range("G5").Value = 128
$ws.Range("H5").Value = '◆ツール,開発'
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5401202')
$ws.Range("F5").Style = "Hyperlink"

# Row 6
$ws.Range("A6").Value = '2025-09-26 12:35:31'
$ws.Range("B6").Value = 'システムの開発補助や運営サポート【フルリモート×長期】'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("F6").Value = 'https://www.lancers.jp/work/detail/5323359'
$ws.Range("G6").Value = 83
$ws.Range("H6").Value = '◆開発'
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5323359')
$ws.Range("F6").Style = "Hyperlink"

# Row 7
$ws.Range("A7").Value = '2025-09-26 12:35:31'
$ws.Range("B7").Value = '【急募】LLMによるMCP(Model Context Protocol)でのExcel操作機能開発'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Range("F7").Value = 'https://www.lancers.jp/work/detail/5400689'
$ws.Range("G7").Value = 75
$ws.Range("H7").Value = '◆開発'
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5400689')
$ws.Range("F7").Style = "Hyperlink"

# Row 8
$ws.Range("A8").Value = '2025-09-26 12:35:31'
$ws.Range("B8").Value = '【急募】Zoho CRMで流入検索キーワード確認設定依頼(zohoコンサルができる方はなお可)'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("F8").Value = 'https://www.lancers.jp/work/detail/5401115'
$ws.Range("G8").Value = 48
$ws.Range("H8").Value = '◆コンサル'
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5401115')
$ws.Range("F8").Style = "Hyperlink"

# Row 9
$ws.Range("A9").Value = '2025-09-26 12:35:31'
$ws.Range("B9").Value = '【急募】音源ライセンス販売サイトのMVP構築依頼'
$ws.Range("C9").Value = 'システム開発'
$ws.Range("D9").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E9").Value = '期限情報なし'
$ws.Range("F9").Value = 'https://www.lancers.jp/work/detail/5400763'
$ws.Range("G9").Value = 45
$ws.Range("H9").Value = '◇サイト'
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5400763')
$ws.Range("F9").Style = "Hyperlink"

# Row 10
$ws.Range("A10").Value = '2025-09-26 12:35:31'
$ws.Range("B10").Value = 'wordpressレンダリングを妨げるリソースの除外'
$ws.Range("C10").Value = 'システム開発'
$ws.Range("D10").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E10").Value = '期限情報なし'
$ws.Range("F10").Value = 'https://www.lancers.jp/work/detail/5016989'
$ws.Range("G10").Value = 33
$ws.Range("H10").Value = '○WordPress'
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.lancers.jp/work/detail/5016989')
$ws.Range("F10").Style = "Hyperlink"

# Row 11
$ws.Range("A11").Value = '2025-09-26 12:35:31'
$ws.Range("B11").Value = 'eBayテラピークでのキーワード検索結果等の取得するためのシステム制作'
$ws.Range("C11").Value = 'システム開発'
$ws.Range("D11").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E11").Value = '期限情報なし'
$ws.Range("F11").Value = 'https://www.lancers.jp/work/detail/5390238'
$ws.Range("G11").Value = 33
$ws.Hyperlinks.Add($ws.Range("F11"), 'https://www.lancers.jp/work/detail/5390238')
$ws.Range("F11").Style = "Hyperlink"

# Row 12
$ws.Range("A12").Value = '2025-09-26 12:35:31'
$ws.Range("B12").Value = 'Drupal関連プロジェクトの要件定義や基本設計ができる方(1人月、長期継続案件)'
$ws.Range("C12").Value = 'システム開発'
$ws.Range("D12").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E12").Value = '期限情報なし'
$ws.Range("F12").Value = 'https://www.lancers.jp/work/detail/5400683'
$ws.Range("G12").Value = 25
$ws.Hyperlinks.Add($ws.Range("F12"), 'https://www.lancers.jp/work/detail/5400683')
$ws.Range("F12").Style = "Hyperlink"

# Row 13
$ws.Range("A13").Value = '2025-09-26 12:35:31'
$ws.Range("B13").Value = '金融関連プロジェクトの要件定義や基本設計ができる方(1人月、長期継続案件)'
$ws.Range("C13").Value = 'システム開発'
$ws.Range("D13").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E13").Value = '期限情報なし'
$ws.Range("F13").Value = 'https://www.lancers.jp/work/detail/5400681'
$ws.Range("G13").Value = 25
$ws.Hyperlinks.Add($ws.Range("F13"), 'https://www.lancers.jp/work/detail/5400681')
$ws.Range("F13").Style = "Hyperlink"

# Row 14
$ws.Range("A14").Value = '2025-09-26 12:35:31'
$ws.Range("B14").Value = '【急募】東京でのWeb制作プロジェクトに参加しませんか?'
$ws.Range("C14").Value = 'システム開発'
$ws.Range("D14").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E14").Value = '期限情報なし'
$ws.Range("F14").Value = 'https://www.lancers.jp/work/detail/5400965'
$ws.Range("G14").Value = 18
$ws.Hyperlinks.Add($ws.Range("F14"), 'https://www.lancers.jp/work/detail/5400965')
$ws.Range("F14").Style = "Hyperlink"

# Row 15
$ws.Range("A15").Value = '2025-09-26 12:35:31'
$ws.Range("B15").Value = '限定公開 PR 限定公開の仕事'
$ws.Range("C15").Value = 'システム開発'
$ws.Range("D15").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E15").Value = '期限情報なし'
$ws.Range("F15").Value = 'https://www.lancers.jp/work/detail/5399347'
$ws.Range("G15").Value = 13
$ws.Hyperlinks.Add($ws.Range("F15"), 'https://www.lancers.jp/work/detail/5399347')
$ws.Range("F15").Style = "Hyperlink"

# Row 16
$ws.Range("A16").Value = '2025-09-26 12:35:31'
$ws.Range("B16").Value = '【急募】スーパードルフィーの洋服をオーダーメイドで作成希望'
$ws.Range("C16").Value = 'システム開発'
$ws.Range("D16").Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Range("E16").Value = '期限情報なし'
$ws.Range("F16").Value = 'https://www.lancers.jp/work/detail/5400988'
$ws.Range("G16").Value = 10
$ws.Hyperlinks.Add($ws.Range("F16"), 'https://www.lancers.jp/work/detail/5400988')
$ws.Range("F16").Style = "Hyperlink"
